{"js": "// Regras de verifica\u00e7\u00e3o e an\u00e1lise de requisitos \u2014 apply the two\n// list-item text rewrites described by the commit.\n//\n// 1) \"Os requisitos devem ser escritos de forma clara e concisa;\" becomes\n//    a longer explanation about precise wording.\n// 2) \"Toda e qualquer mudan\u00e7a em qualquer um dos requisitos deve ser\n//    documentada;\" keeps the same wording but is retyped as a single run\n//    (clearing the stray proofing-error markers around \"documentada\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst OLD_TEXT_1 = \"Os requisitos devem ser escritos de forma clara e concisa;\";\nconst NEW_TEXT_1 =\n  \"Os requisitos devem ser escritos de forma precisa, isto \u00e9, evitar \" +\n  \"palavras que possam gerar mais de uma interpreta\u00e7\u00e3o, tais como \" +\n  \"\\u201ce\\u201d, \\u201cou\\u201d, \\u201cmas\\u201d, \\\"amig\u00e1vel\\u201d, \" +\n  \"\\u201cflex\u00edvel\\u201d entre outras similares.\";\n\nconst OLD_TEXT_2 =\n  \"Toda e qualquer mudan\u00e7a em qualquer um dos requisitos deve ser documentada;\";\nconst NEW_TEXT_2 = OLD_TEXT_2; // text unchanged, only the run/markup is normalized\n\nlet paragraph1 = null;\nlet paragraph2 = null;\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === OLD_TEXT_1) {\n    paragraph1 = paragraph;\n  } else if (paragraph.text === OLD_TEXT_2) {\n    paragraph2 = paragraph;\n  }\n}\n\nif (!paragraph1) {\n  throw new Error('Paragraph \"' + OLD_TEXT_1 + '\" not found.');\n}\nif (!paragraph2) {\n  throw new Error('Paragraph \"' + OLD_TEXT_2 + '\" not found.');\n}\n\n// \"Replace\" rewrites the paragraph's range as fresh run(s), which is what\n// naturally drops the old <w:proofErr/> bookkeeping around \"documentada\".\nparagraph1.insertText(NEW_TEXT_1, \"Replace\");\nparagraph2.insertText(NEW_TEXT_2, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Regras de verifica\u00e7\u00e3o e an\u00e1lise de requisitos \u2014 apply the two\n# list-item text rewrites described by the commit.\n#\n# 1) \"Os requisitos devem ser escritos de forma clara e concisa;\" becomes\n#    a longer explanation about precise wording.\n# 2) \"Toda e qualquer mudan\u00e7a em qualquer um dos requisitos deve ser\n#    documentada;\" keeps the same wording but is retyped as a single run\n#    (clearing the stray proofing-error markers around \"documentada\").\n\n$d = $word.ActiveDocument\n\n$oldText1 = \"Os requisitos devem ser escritos de forma clara e concisa;\"\n$newText1 = 'Os requisitos devem ser escritos de forma precisa, isto \u00e9, evitar palavras que possam gerar mais de uma interpreta\u00e7\u00e3o, tais como \u201ce\u201d, \u201cou\u201d, \u201cmas\u201d, \"amig\u00e1vel\u201d, \u201cflex\u00edvel\u201d entre outras similares.'\n\n$oldText2 = \"Toda e qualquer mudan\u00e7a em qualquer um dos requisitos deve ser documentada;\"\n$newText2 = \"Toda e qualquer mudan\u00e7a em qualquer um dos requisitos deve ser documentada;\"\n\n$paragraph1 = $null\n$paragraph2 = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $oldText1) {\n        $paragraph1 = $p\n    } elseif ($t -eq $oldText2) {\n        $paragraph2 = $p\n    }\n}\n\nif ($paragraph1 -eq $null) {\n    throw \"Paragraph '$oldText1' not found.\"\n}\nif ($paragraph2 -eq $null) {\n    throw \"Paragraph '$oldText2' not found.\"\n}\n\n# Shrink the range so it excludes the trailing paragraph mark, clear it, then\n# retype the full sentence. Retyping (rather than a straight Text\n# assignment) is what makes Word rebuild the paragraph as a single fresh\n# run, dropping the old <w:proofErr/> bookkeeping around \"documentada\".\n$range1 = $paragraph1.Range\n$range1.MoveEnd(1, -1) | Out-Null\n$range1.Text = \"\"\n$range1.Text = $newText1\n\n$range2 = $paragraph2.Range\n$range2.MoveEnd(1, -1) | Out-Null\n$range2.Text = \"\"\n$range2.Text = $newText2\n"}
